$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: fill in the journal entry for 13.03.20 (CLion / Bataille Navale / Programmation du jeu)
$ws.Range("B10").Value = 43903
$ws.Range("C10").Value = 0.40972222222222227
$ws.Range("D10").Value = 0.51041666666666663
$ws.Range("F10").Value = "CLion"
$ws.Range("G10").Value = "Bataille Navale"
$ws.Range("H10").Value = "Programmation du jeu"
$ws.Range("I10").Value = "Création des différentes fonctions"

# Move the active selection to I11
$ws.Range("I11").Select()
